$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing data row (row 2) down into the new
# row 3 so the new cells pick up the same cell style (s="1") as the rest
# of the table.
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)

# Populate the new record: "arroz integral" / comida / carboidrato with
# its nutritional values.
$ws.Range("A3").Value = "arroz integral"
$ws.Range("B3").Value = "comida"
$ws.Range("C3").Value = "carboidrato"
$ws.Range("D3").Value = 124
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2.6
$ws.Range("G3").Value = 25.8
